$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure column D values (which look numeric, e.g. "43.617.76" or "1.01")
# are stored as text, matching the source data which uses inline strings
# rather than numbers (prices here use "." as a thousands separator).
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "43.617.76"
$ws.Range("E2").Value = "  +1.26%  "
$ws.Range("D3").Value = "2.270.32"
$ws.Range("E3").Value = "  -0.13%  "
$ws.Range("E4").Value = "  +0.11%  "
$ws.Range("D5").Value = "117.79"
$ws.Range("E5").Value = "  +6.01%  "
$ws.Range("D6").Value = "267.84"
$ws.Range("E6").Value = "  +1.05%  "
$ws.Range("D7").Value = "0.641"
$ws.Range("E7").Value = "  +3.56%  "
$ws.Range("D8").Value = "1.01"
$ws.Range("E8").Value = "  +0.31%  "
$ws.Range("D9").Value = "0.620"
$ws.Range("E9").Value = "  +2.58%  "
$ws.Range("D10").Value = "47.30"
$ws.Range("E10").Value = "  -0.46%  "
$ws.Range("D11").Value = "0.0941"
$ws.Range("E11").Value = "  +1.29%  "
$ws.Range("D12").Value = "9.40"
$ws.Range("E12").Value = "  +7.30%  "
$ws.Range("E13").Value = "  -1.46%  "
$ws.Range("D14").Value = "15.68"
$ws.Range("E14").Value = "  +1.59%  "
$ws.Range("D15").Value = "0.903"
$ws.Range("E15").Value = "  +5.92%  "
$ws.Range("D16").Value = "2.616.82"
$ws.Range("E16").Value = "  +0.10%  "
$ws.Range("D17").Value = "2.270.69"
$ws.Range("E17").Value = "  -0.28%  "
$ws.Range("D18").Value = "43.613.43"
$ws.Range("E18").Value = "  +1.18%  "
$ws.Range("D19").Value = "0.0000109"
$ws.Range("E19").Value = "  +1.43%  "
$ws.Range("D20").Value = "6.89"
$ws.Range("E20").Value = "  +1.08%  "
$ws.Range("E21").Value = "  +1.82%  "
$ws.Range("D22").Value = "2.39"
$ws.Range("E22").Value = "  -4.94%  "
$ws.Range("D23").Value = "234.30"
$ws.Range("E23").Value = "  +1.25%  "
$ws.Range("E24").Value = "  +2.82%  "
$ws.Range("D25").Value = "9.67"
$ws.Range("E25").Value = "  +0.17%  "
$ws.Range("D26").Value = "12.21"
$ws.Range("E26").Value = "  +8.10%  "
$ws.Range("E27").Value = "  +1.88%  "
$ws.Range("D28").Value = "41.62"
$ws.Range("E28").Value = "  +2.86%  "
$ws.Range("D29").Value = "3.35"
$ws.Range("E29").Value = "  +2.37%  "
$ws.Range("D31").Value = "174.17"
$ws.Range("E31").Value = "  +1.45%  "
$ws.Range("D32").Value = "21.49"
$ws.Range("E32").Value = "  +0.78%  "
$ws.Range("D33").Value = "0.0918"
$ws.Range("E33").Value = "  +1.88%  "
$ws.Range("D34").Value = "5.72"
$ws.Range("E34").Value = "  +0.38%  "
$ws.Range("E35").Value = "  +2.85%  "
$ws.Range("D36").Value = "4.28"
$ws.Range("E36").Value = "  +12.71%  "
$ws.Range("D37").Value = "0.0383"
$ws.Range("E37").Value = "  +9.06%  "
$ws.Range("E38").Value = "  -1.55%  "
$ws.Range("E39").Value = "  +3.07%  "
$ws.Range("E40").Value = "  -2.46%  "
$ws.Range("D41").Value = "13.81"
$ws.Range("E41").Value = "  -0.20%  "
$ws.Range("D42").Value = "0.239"
$ws.Range("E42").Value = "  +1.64%  "
$ws.Range("D43").Value = "72.02"
$ws.Range("E43").Value = "  -5.54%  "
$ws.Range("E44").Value = "  +0.15%  "
$ws.Range("E45").Value = "  +1.00%  "
$ws.Range("D46").Value = "5.74"
$ws.Range("E46").Value = "  -5.50%  "
$ws.Range("D47").Value = "0.673"
$ws.Range("E47").Value = "  +20.57%  "
$ws.Range("D48").Value = "73.90"
$ws.Range("E48").Value = "  +38.69%  "
$ws.Range("E49").Value = "  +1.83%  "
$ws.Range("B50").Value = "Aave"
$ws.Range("C50").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D50").Value = "103.02"
$ws.Range("E50").Value = "  +1.95%  "
$ws.Range("B51").Value = "FraxShare"
$ws.Range("C51").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D51").Value = "8.56"
$ws.Range("E51").Value = "  -0.78%  "
